$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor values updated
$ws.Range("B3").Value = 0.998881982275008
$ws.Range("C3").Value = 0.9988402941198884
$ws.Range("D3").Value = 0.997298858750292

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9987236165483336
$ws.Range("C4").Value = 0.9987068685116811
$ws.Range("D4").Value = 0.9959935393512063

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9991083750850099
$ws.Range("C5").Value = 0.9991877088299163
$ws.Range("D5").Value = 0.9991058658974094
